# Updated mapping of SubjectID and RecordID
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# The NIEM 3.1 Mapping XPaths for "Subject ID" (row 15) and
# "Booking Number" / agency record id (row 16) were pointing at the wrong
# root document element (br-doc:BookingReport). Point them at the
# Consent Decision Report document instead.
$ws.Range("C15").Value = "/cdr-doc:ConsentDecisionReport/j:Booking/j:BookingSubject/j:SubjectIdentification/nc:IdentificationID"
$ws.Range("C16").Value = "cdr-doc:ConsentDecisionReport/j:Booking/j:BookingAgencyRecordIdentification/nc:IdentificationID"

# Reflect the new focus cell / window position as left by the editor.
$ws.Range("C16").Select() | Out-Null
